$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.877.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.770.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4480"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07451"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.094"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.045"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.220"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.773.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06419"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.821"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.896.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.116"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.977.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.192"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.099"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09157"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.576"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02290"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06095"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2097"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6330"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.959"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.392"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.924"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.739"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5874"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.959"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06915"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.139"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.76%  "
